$d = $word.ActiveDocument

# The document contains placeholder runs split as "<id>" / "p160r_1" / "</id>"
# (and similarly for "p160r_2"). Collapse each triplet into a single run
# "<id>p160r_1</id>" by doing a self-replace via Find/Execute, which causes
# Word to merge the matched text into one run using the first run's
# formatting. We must NOT touch the "fig_p160r_1" / "fig_p160r_2" variants,
# so the search text is anchored exactly to "<id>p160r_1</id>" and
# "<id>p160r_2</id>".

$targets = @("p160r_1", "p160r_2")

foreach ($t in $targets) {
    $search = "<id>" + $t + "</id>"
    $d.Content.Find.Execute($search, $false, $false, $false, $false, $false, $true, 1, $false, $search, 2) | Out-Null
}
